$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string rich text edits (Volume/date header) ---
$a8 = $ws.Range("A8")
$a8.Characters(21, 2).Text = "49"

$c9 = $ws.Range("C9")
$c9.Characters(27, 10).Text = "12/5/2022"
$c9.Characters(47, 9).Text = "12/11/2022"

# --- Helper to convert a numeric cell into a shared-text cell ("0" or "***.*") ---
function Set-TextCell($cellRef, $sourceRef) {
    $ws.Range($sourceRef).Copy()
    $ws.Range($cellRef).PasteSpecial(-4122)
    $ws.Range($sourceRef).Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
}

# --- Cells that change from numeric to shared text ---
Set-TextCell "D22" "C26"
Set-TextCell "E22" "E26"
Set-TextCell "C27" "C26"
Set-TextCell "D27" "C26"
Set-TextCell "E27" "E26"
Set-TextCell "D28" "C26"
Set-TextCell "E28" "E26"
Set-TextCell "D29" "C26"
Set-TextCell "E29" "E26"

# --- Numeric value updates ---
$ws.Range("M15").Value = -37.037037037037
$ws.Range("N15").Value = -72.131147540983
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = -83.333333333333
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 19
$ws.Range("H16").Value = -47.368421052631
$ws.Range("I16").Value = 181
$ws.Range("J16").Value = 179
$ws.Range("K16").Value = 1.117318435754
$ws.Range("L16").Value = 15.286624203821
$ws.Range("M16").Value = -25.819672131147
$ws.Range("N16").Value = -79.266895761741
$ws.Range("C17").Value = 9
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = 50
$ws.Range("F17").Value = 34
$ws.Range("G17").Value = 29
$ws.Range("H17").Value = 17.241379310344
$ws.Range("I17").Value = 372
$ws.Range("J17").Value = 352
$ws.Range("K17").Value = 5.681818181818
$ws.Range("L17").Value = 26.530612244898
$ws.Range("M17").Value = 73.023255813953
$ws.Range("N17").Value = -33.093525179856
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 25
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = 108.333333333333
$ws.Range("I18").Value = 194
$ws.Range("J18").Value = 103
$ws.Range("K18").Value = 88.349514563106
$ws.Range("L18").Value = 31.972789115646
$ws.Range("M18").Value = 120.454545454545
$ws.Range("N18").Value = -57.362637362637
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 30
$ws.Range("G19").Value = 38
$ws.Range("H19").Value = -21.052631578947
$ws.Range("I19").Value = 459
$ws.Range("J19").Value = 355
$ws.Range("K19").Value = 29.295774647887
$ws.Range("L19").Value = 37.837837837837
$ws.Range("M19").Value = 87.34693877551
$ws.Range("N19").Value = -15.469613259668
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = 60
$ws.Range("I20").Value = 82
$ws.Range("J20").Value = 58
$ws.Range("K20").Value = 41.379310344827
$ws.Range("L20").Value = 28.125
$ws.Range("M20").Value = 121.621621621622
$ws.Range("N20").Value = -79.448621553884
$ws.Range("C21").Value = 24
$ws.Range("D21").Value = 29
$ws.Range("E21").Value = -17.241379310344
$ws.Range("F21").Value = 107
$ws.Range("G21").Value = 105
$ws.Range("H21").Value = 1.904761904761
$ws.Range("I21").Value = 1309
$ws.Range("J21").Value = 1071
$ws.Range("K21").Value = 22.222222222222
$ws.Range("L21").Value = 28.333333333333
$ws.Range("M21").Value = 52.386495925494
$ws.Range("N21").Value = -55.063508410573
$ws.Range("F22").Value = 2
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 8
$ws.Range("K22").Value = 14.285714285714
$ws.Range("L22").Value = -11.111111111111
$ws.Range("M22").Value = 33.333333333333
$ws.Range("C23").Value = 5
$ws.Range("D23").Value = 8
$ws.Range("E23").Value = -37.5
$ws.Range("G23").Value = 39
$ws.Range("H23").Value = -28.205128205128
$ws.Range("I23").Value = 388
$ws.Range("J23").Value = 412
$ws.Range("K23").Value = -5.825242718446
$ws.Range("L23").Value = -3.722084367245
$ws.Range("M23").Value = 39.068100358422
$ws.Range("C24").Value = 27
$ws.Range("D24").Value = 14
$ws.Range("E24").Value = 92.857142857142
$ws.Range("F24").Value = 96
$ws.Range("G24").Value = 51
$ws.Range("H24").Value = 88.235294117647
$ws.Range("I24").Value = 877
$ws.Range("J24").Value = 717
$ws.Range("K24").Value = 22.31520223152
$ws.Range("L24").Value = 21.300138312586
$ws.Range("M24").Value = 37.245696400626
$ws.Range("C25").Value = 11
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = 57.142857142857
$ws.Range("F25").Value = 34
$ws.Range("G25").Value = 42
$ws.Range("H25").Value = -19.047619047619
$ws.Range("I25").Value = 516
$ws.Range("J25").Value = 513
$ws.Range("K25").Value = 0.584795321637
$ws.Range("L25").Value = -1.149425287356
$ws.Range("M25").Value = -18.867924528301
$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 9
$ws.Range("H27").Value = -44.444444444444
$ws.Range("L27").Value = 34.042553191489
$ws.Range("F28").Value = 2
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = -50
$ws.Range("I28").Value = 26
$ws.Range("K28").Value = -23.529411764705
$ws.Range("L28").Value = -36.585365853658
$ws.Range("M28").Value = -33.333333333333
$ws.Range("N28").Value = -72.340425531914
$ws.Range("F29").Value = 2
$ws.Range("G29").Value = 4
$ws.Range("H29").Value = -50
$ws.Range("I29").Value = 19
$ws.Range("K29").Value = -40.625
$ws.Range("L29").Value = -44.117647058823
$ws.Range("M29").Value = -44.117647058823
$ws.Range("N29").Value = -78.40909090909
